$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: header "P_E" ---
$ws.Range("E1").Value = "P_E"
# Copy the header formatting (bold/border/alignment) from D1 onto E1
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 2 data refresh ---
$ws.Range("A2").Value = "Coca-Cola Company The"

# Force plain-text storage (leading apostrophe = Excel's "treat as text" marker)
# then reset the cell style so no stray number-format / quote-prefix sticks around.
$ws.Range("B2").Value = "'$70.00"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "-0.02(0.03%) 1D"

# D2 stays "N/A" (unchanged)

$ws.Range("E2").Value = "'28.3"
$ws.Range("E2").Style = "Normal"
